$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 23666.666
$ws.Range("I21").Value = 20000
$ws.Range("K21").Value = 20000
$ws.Range("M21").Value = -19532
$ws.Range("H23").Value = 23666.666
$ws.Range("I23").Value = 20000
$ws.Range("K23").Value = 20000
$ws.Range("M23").Value = -19766
$ws.Range("H29").Value = 775
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 3000
$ws.Range("N29").Value = -3562
$ws.Range("H32").Value = 983.1667
$ws.Range("I32").Value = 633
$ws.Range("K32").Value = 633
$ws.Range("M32").Value = -307
$ws.Range("H38").Value = 143.27272
$ws.Range("I38").Value = 143.27272
$ws.Range("K38").Value = 429.81816
$ws.Range("M38").Value = -57.81815999999998
$ws.Range("H53").Value = 332.91666
$ws.Range("I53").Value = 302.14285
$ws.Range("K53").Value = 302.14285
$ws.Range("M53").Value = 334.85715
$ws.Range("H58").Value = 596.5789
$ws.Range("I58").Value = 74.72221999999999
$ws.Range("J58").Value = 9990
$ws.Range("K58").Value = 224.16666
$ws.Range("L58").Value = 29970
$ws.Range("M58").Value = -74.16665999999998
$ws.Range("N58").Value = -30270
$ws.Range("H106").Value = 4000
$ws.Range("I106").Value = 4000
$ws.Range("K106").Value = 4000
$ws.Range("M106").Value = -3369
$ws.Range("H132").Value = 29616.621
$ws.Range("I132").Value = 34923.13
$ws.Range("K132").Value = 104769.39
$ws.Range("M132").Value = -102239.39
$ws.Range("H137").Value = 1004.8485
$ws.Range("I137").Value = 721.53845
$ws.Range("J137").Value = 2057.1428
$ws.Range("K137").Value = 2164.61535
$ws.Range("L137").Value = 6171.428400000001
$ws.Range("M137").Value = 385.38465
$ws.Range("N137").Value = -11271.4284
$ws.Range("H138").Value = 4421.7896
$ws.Range("I138").Value = 10600
$ws.Range("J138").Value = 3694.9412
$ws.Range("K138").Value = 31800
$ws.Range("L138").Value = 11084.8236
$ws.Range("M138").Value = -26660
$ws.Range("N138").Value = -21364.8236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2045.4
$ws.Range("I32").Value = 1350.6774
$ws.Range("K32").Value = 1350.6774
$ws.Range("M32").Value = -1063.6774
$ws.Range("H45").Value = 1981.7838
$ws.Range("I45").Value = 1504.25
$ws.Range("J45").Value = 2863.3845
$ws.Range("K45").Value = 1504.25
$ws.Range("L45").Value = 2863.3845
$ws.Range("M45").Value = -1127.25
$ws.Range("N45").Value = -3617.3845
$ws.Range("H132").Value = 15974.771
$ws.Range("I132").Value = 1193.24
$ws.Range("J132").Value = 52928.6
$ws.Range("K132").Value = 3579.72
$ws.Range("L132").Value = 158785.8
$ws.Range("M132").Value = -1049.72
$ws.Range("N132").Value = -163845.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1157.6666
$ws.Range("J16").Value = 1474.75
$ws.Range("L16").Value = 1474.75
$ws.Range("N16").Value = -2048.75
$ws.Range("H22").Value = 537.1111
$ws.Range("I22").Value = 540.6667
$ws.Range("J22").Value = 530
$ws.Range("K22").Value = 540.6667
$ws.Range("L22").Value = 530
$ws.Range("M22").Value = -190.6667
$ws.Range("N22").Value = -1230
$ws.Range("H35").Value = 854
$ws.Range("I35").Value = 854
$ws.Range("K35").Value = 854
$ws.Range("M35").Value = -560
$ws.Range("H113").Value = 1157.6666
$ws.Range("J113").Value = 1474.75
$ws.Range("L113").Value = 1474.75
$ws.Range("N113").Value = -5814.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 241.83333
$ws.Range("I40").Value = 90.2
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 360.8
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -291.8
$ws.Range("N40").Value = -4138
$ws.Range("H131").Value = 806.11224
$ws.Range("I131").Value = 246
$ws.Range("J131").Value = 823.8
$ws.Range("K131").Value = 738
$ws.Range("L131").Value = 2471.4
$ws.Range("M131").Value = 4302
$ws.Range("N131").Value = -12551.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2154.8096
$ws.Range("I122").Value = 1697.421
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 5092.263
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -2642.263
$ws.Range("N122").Value = -24400
$ws.Range("H132").Value = 19998.7
$ws.Range("I132").Value = 3825.5789
$ws.Range("J132").Value = 47934.09
$ws.Range("K132").Value = 11476.7367
$ws.Range("L132").Value = 143802.27
$ws.Range("M132").Value = -8946.736699999999
$ws.Range("N132").Value = -148862.27

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4765.9375
$ws.Range("I7").Value = 4734.615
$ws.Range("J7").Value = 4901.6665
$ws.Range("K7").Value = 4734.615
$ws.Range("L7").Value = 4901.6665
$ws.Range("M7").Value = -4622.615
$ws.Range("N7").Value = -5125.6665
$ws.Range("H16").Value = 1132.8572
$ws.Range("H61").Value = 4379
$ws.Range("I61").Value = 2482.3125
$ws.Range("J61").Value = 8714.286
$ws.Range("K61").Value = 2482.3125
$ws.Range("L61").Value = 8714.286
$ws.Range("M61").Value = -2280.3125
$ws.Range("N61").Value = -9118.286
$ws.Range("H113").Value = 4379
$ws.Range("I113").Value = 2482.3125
$ws.Range("J113").Value = 8714.286
$ws.Range("K113").Value = 2482.3125
$ws.Range("L113").Value = 8714.286
$ws.Range("M113").Value = -312.3125
$ws.Range("N113").Value = -13054.286
$ws.Range("H126").Value = 4765.9375
$ws.Range("I126").Value = 4734.615
$ws.Range("J126").Value = 4901.6665
$ws.Range("K126").Value = 14203.845
$ws.Range("L126").Value = 14704.9995
$ws.Range("M126").Value = -11733.845
$ws.Range("N126").Value = -19644.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3004149.8
$ws.Range("I113").Value = 1928.4
$ws.Range("J113").Value = 6756926.5
$ws.Range("K113").Value = 5785.200000000001
$ws.Range("L113").Value = 20270779.5
$ws.Range("M113").Value = -3615.200000000001
$ws.Range("N113").Value = -20275119.5
$ws.Range("H136").Value = 1793734.2
$ws.Range("I136").Value = 4033246
$ws.Range("J136").Value = 2124.9
$ws.Range("K136").Value = 12099738
$ws.Range("L136").Value = 6374.700000000001
$ws.Range("M136").Value = -12097188
$ws.Range("N136").Value = -11474.7
